$d = $word.ActiveDocument

# The document currently has a single introductory paragraph. Append five
# new bulleted list paragraphs -- one per file in the portfolio project
# folder -- each formatted with the "List Paragraph" style and a bullet
# numbering definition, matching how Word's default bullet-list gallery
# button behaves.

$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Project_"
$p2.Range.InsertAfter("Table_Creation.sql")
$p2.Range.ListFormat.ApplyBulletDefault()
Write-Host "list item 1 done"
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Project_Data_Insert"
$p3.Range.InsertAfter(".sql")
$p3.Range.ListFormat.ApplyBulletDefault()
Write-Host "list item 2 done"
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$p4.Range.Text = "Project_Views.sql"
$p4.Range.ListFormat.ApplyBulletDefault()
Write-Host "list item 3 done"
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5.Range.Text = "Project_Database.accdb"
$p5.Range.ListFormat.ApplyBulletDefault()
Write-Host "list item 4 done"
$p5.Range.InsertParagraphAfter()

$p6 = $d.Paragraphs(6)
$p6.Range.Text = "Seshu_Miriyala_Project.docx"
$p6.Range.ListFormat.ApplyBulletDefault()
Write-Host "list item 5 done"
